$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '67.999.78'
$ws.Range("E2").Value = '  +0.49%  '
$ws.Range("D3").Value = '3.304.04'
$ws.Range("E3").Value = '  +1.14%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("E4").Value = '  -0.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '186.47'
$ws.Range("E5").Value = '  +1.22%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '582.57'
$ws.Range("E6").Value = '  +0.31%  '
$ws.Range("E7").Value = '  +0.01%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.597'
$ws.Range("E8").Value = '  -1.06%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.130'
$ws.Range("E9").Value = '  +0.00%  '
$ws.Range("E10").Value = '  +1.40%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.408'
$ws.Range("E11").Value = '  +0.00%  '
$ws.Range("D12").Value = '3.878.33'
$ws.Range("E12").Value = '  +1.24%  '
$ws.Range("E13").Value = '  -2.13%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '27.51'
$ws.Range("E14").Value = '  +0.68%  '
$ws.Range("D15").Value = '68.137.35'
$ws.Range("E15").Value = '  +0.60%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000168'
$ws.Range("E16").Value = '  +0.07%  '
$ws.Range("D17").Value = '3.307.64'
$ws.Range("E17").Value = '  +1.64%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '447.63'
$ws.Range("E18").Value = '  +12.06%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '5.71'
$ws.Range("E19").Value = '  -0.09%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '13.52'
$ws.Range("E20").Value = '  +0.66%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '7.74'
$ws.Range("E21").Value = '  +2.41%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '74.93'
$ws.Range("E22").Value = '  +5.72%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("E23").Value = '  +0.04%  '
$ws.Range("D24").Value = '3.462.46'
$ws.Range("E24").Value = '  +1.67%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '0.514'
$ws.Range("E25").Value = '  +1.26%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.0000119'
$ws.Range("E26").Value = '  +1.44%  '
$ws.Range("E27").Value = '  +0.35%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '9.12'
$ws.Range("E28").Value = '  -4.36%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.00'
$ws.Range("E29").Value = '  +0.10%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '1.98'
$ws.Range("E30").Value = '  +1.50%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '22.86'
$ws.Range("E31").Value = '  +0.99%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '5.35'
$ws.Range("E32").Value = '  -2.05%  '
$ws.Range("E33").Value = '  +0.02%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.25'
$ws.Range("E34").Value = '  +0.45%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.80'
$ws.Range("E35").Value = '  -1.87%  '
$ws.Range("E36").Value = '  +4.70%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '163.55'
$ws.Range("E37").Value = '  -0.01%  '
$ws.Range("E38").Value = '  -1.12%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '27.01'
$ws.Range("E39").Value = '  +0.17%  '
$ws.Range("E40").Value = '  +0.02%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.782'
$ws.Range("E41").Value = '  -2.84%  '
$ws.Range("E42").Value = '  +2.13%  '
$ws.Range("D43").Value = '2.703.32'
$ws.Range("E43").Value = '  +1.10%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '40.71'
$ws.Range("E44").Value = '  +0.02%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0673'
$ws.Range("E45").Value = '  -0.85%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.41'
$ws.Range("E46").Value = '  -0.83%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '24.68'
$ws.Range("E47").Value = '  +0.71%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '327.02'
$ws.Range("E48").Value = '  -2.52%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '0.0276'
$ws.Range("E49").Value = '  +1.07%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '31.57'
$ws.Range("E50").Value = '  +3.28%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.989'
$ws.Range("E51").Value = '  +2.16%  '
